$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Unit Test Case")
$ws.Activate()

# Fix: D32 value changes from "N" to "Y"
$ws.Range("D32").Value = "Y"

# Append new API test rows (70-99) to the Unit Test Case sheet
$ws.Range("A70").Value = "Get snaps by snap product id"
$ws.Range("B70").Value = 200
$ws.Range("C70").Value = "Normal"
$ws.Range("D70").Value = "Y"

$ws.Range("B71").Value = 404
$ws.Range("C71").Value = "Unexisting snap_product_id"
$ws.Range("D71").Value = "Y"

$ws.Range("B72").Value = 404
$ws.Range("C72").Value = "Missing snap_product_id"
$ws.Range("D72").Value = "Y"

$ws.Range("A74").Value = "Get user profile"
$ws.Range("B74").Value = 200
$ws.Range("C74").Value = "Self user_id"
$ws.Range("D74").Value = "Y"

$ws.Range("B75").Value = 200
$ws.Range("C75").Value = "Another valid user_id"
$ws.Range("D75").Value = "Y"

$ws.Range("B76").Value = 404
$ws.Range("C76").Value = "Unexisting user_id"
$ws.Range("D76").Value = "Y"

$ws.Range("B77").Value = 404
$ws.Range("C77").Value = "Missing user_id"
$ws.Range("D77").Value = "Y"

$ws.Range("A79").Value = "Change password"
$ws.Range("B79").Value = 200
$ws.Range("C79").Value = "Normal"
$ws.Range("D79").Value = "Y"

$ws.Range("B80").Value = 400
$ws.Range("C80").Value = "Missing current password"
$ws.Range("D80").Value = "Y"

$ws.Range("B81").Value = 400
$ws.Range("C81").Value = "Missing new password"
$ws.Range("D81").Value = "Y"

$ws.Range("B82").Value = 400
$ws.Range("C82").Value = "Invalid new password"
$ws.Range("D82").Value = "Y"

$ws.Range("B83").Value = 401
$ws.Range("C83").Value = "Unauthorized user_id"
$ws.Range("D83").Value = "Y"

$ws.Range("B84").Value = 401
$ws.Range("C84").Value = "Logout and Change password"
$ws.Range("D84").Value = "Y"

$ws.Range("B85").Value = 404
$ws.Range("C85").Value = "Invalid current password"
$ws.Range("D85").Value = "Y"

$ws.Range("A87").Value = "Update user profile"
$ws.Range("B87").Value = 200
$ws.Range("C87").Value = "Firstname = 'Testing'"
$ws.Range("D87").Value = "Y"

$ws.Range("B88").Value = 200
$ws.Range("C88").Value = "Missing username"
$ws.Range("D88").Value = "Y"

$ws.Range("B89").Value = 200
$ws.Range("C89").Value = "Invalid username"
$ws.Range("D89").Value = "Y"

$ws.Range("B90").Value = 200
$ws.Range("C90").Value = "over 300 character bio"
$ws.Range("D90").Value = "Y"

$ws.Range("B91").Value = 200
$ws.Range("C91").Value = "Logout and update"
$ws.Range("D91").Value = "Y"

$ws.Range("B92").Value = 200
$ws.Range("C92").Value = "Unauthorized user_id"
$ws.Range("D92").Value = "Y"

$ws.Range("B93").Value = 200
$ws.Range("C93").Value = "Missing user_id"
$ws.Range("D93").Value = "Y"

$ws.Range("B94").Value = 200
$ws.Range("C94").Value = "Invalid user_id"
$ws.Range("D94").Value = "Y"

$ws.Range("A96").Value = "Update user profile picture"
$ws.Range("B96").Value = 200
$ws.Range("C96").Value = "10kB image"
$ws.Range("D96").Value = "Y"

$ws.Range("B97").Value = 400
$ws.Range("C97").Value = "Missing image name"
$ws.Range("D97").Value = "Y"

$ws.Range("B98").Value = 400
$ws.Range("C98").Value = "Missing image body"
$ws.Range("D98").Value = "Y"

$ws.Range("B99").Value = 400
$ws.Range("C99").Value = "Invalid image body"
$ws.Range("D99").Value = "Y"

# Update the active selection/scroll position to match the saved view state
$win = $excel.ActiveWindow
$win.ScrollRow = 85
$win.ScrollColumn = 1
$ws.Range("C97").Select()
